$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Section_A")
$ws2 = $wb.Worksheets.Item("Section_B")

# --- Section_A (sheet1) ---
$ws1.Range("A2").Value = '09:00-10:30'
$ws1.Range("B2").Value = 'Free'
$ws1.Range("C2").Value = 'CS263'
$ws1.Range("D2").Value = 'Free'
$ws1.Range("E2").Value = 'CS263'
$ws1.Range("F2").Value = 'CS264'
$ws1.Range("A3").Value = '10:30-12:00'
$ws1.Range("B3").Value = 'CS263'
$ws1.Range("C3").Value = 'Free'
$ws1.Range("D3").Value = 'MA261'
$ws1.Range("E3").Value = 'Free'
$ws1.Range("F3").Value = 'MA261'
$ws1.Range("A4").Value = '12:00-13:00'
$ws1.Range("B4").Value = 'LUNCH BREAK'
$ws1.Range("C4").Value = 'LUNCH BREAK'
$ws1.Range("D4").Value = 'LUNCH BREAK'
$ws1.Range("E4").Value = 'LUNCH BREAK'
$ws1.Range("F4").Value = 'LUNCH BREAK'
$ws1.Range("A5").Value = '13:00-14:30'
$ws1.Range("B5").Value = 'CS261'
$ws1.Range("C5").Value = 'Free'
$ws1.Range("D5").Value = 'CS264'
$ws1.Range("E5").Value = 'CS264'
$ws1.Range("F5").Value = 'Free'
$ws1.Range("A6").Value = '14:30-15:30'
$ws1.Range("B6").Value = 'Free'
$ws1.Range("C6").Value = 'Free'
$ws1.Range("D6").Value = 'Free'
$ws1.Range("E6").Value = 'Free'
$ws1.Range("F6").Value = 'Free'
$ws1.Range("A7").Value = '15:30-17:00'
$ws1.Range("B7").Value = 'Free'
$ws1.Range("C7").Value = 'CS261'
$ws1.Range("D7").Value = 'Free'
$ws1.Range("E7").Value = 'Free'
$ws1.Range("F7").Value = 'CS261'
$ws1.Range("A8").Value = '17:00-18:00'
$ws1.Range("B8").Value = 'Free'
$ws1.Range("C8").Value = 'Free'
$ws1.Range("D8").Value = 'Free'
$ws1.Range("E8").Value = 'Free'
$ws1.Range("F8").Value = 'CS264 (Tutorial)'

# --- Section_B (sheet2) ---
$ws2.Range("A2").Value = '09:00-10:30'
$ws2.Range("B2").Value = 'Free'
$ws2.Range("C2").Value = 'Free'
$ws2.Range("D2").Value = 'Free'
$ws2.Range("E2").Value = 'Free'
$ws2.Range("F2").Value = 'Free'
$ws2.Range("A3").Value = '10:30-12:00'
$ws2.Range("B3").Value = 'CS261'
$ws2.Range("C3").Value = 'CS263'
$ws2.Range("D3").Value = 'Free'
$ws2.Range("E3").Value = 'CS263'
$ws2.Range("F3").Value = 'MA261'
$ws2.Range("A4").Value = '12:00-13:00'
$ws2.Range("B4").Value = 'LUNCH BREAK'
$ws2.Range("C4").Value = 'LUNCH BREAK'
$ws2.Range("D4").Value = 'LUNCH BREAK'
$ws2.Range("E4").Value = 'LUNCH BREAK'
$ws2.Range("F4").Value = 'LUNCH BREAK'
$ws2.Range("A5").Value = '13:00-14:30'
$ws2.Range("B5").Value = 'CS263'
$ws2.Range("C5").Value = 'Free'
$ws2.Range("D5").Value = 'Free'
$ws2.Range("E5").Value = 'CS261'
$ws2.Range("F5").Value = 'CS261'
$ws2.Range("A6").Value = '14:30-15:30'
$ws2.Range("B6").Value = 'CS264 (Tutorial)'
$ws2.Range("C6").Value = 'Free'
$ws2.Range("D6").Value = 'Free'
$ws2.Range("E6").Value = 'Free'
$ws2.Range("F6").Value = 'Free'
$ws2.Range("A7").Value = '15:30-17:00'
$ws2.Range("B7").Value = 'CS264'
$ws2.Range("C7").Value = 'Free'
$ws2.Range("D7").Value = 'CS264'
$ws2.Range("E7").Value = 'MA261'
$ws2.Range("F7").Value = 'CS264'
$ws2.Range("A8").Value = '17:00-18:00'
$ws2.Range("B8").Value = 'Free'
$ws2.Range("C8").Value = 'Free'
$ws2.Range("D8").Value = 'Free'
$ws2.Range("E8").Value = 'Free'
$ws2.Range("F8").Value = 'Free'

# Clear old rows 9-12 (previously used, now removed) and update dimension
$ws1.Range("A9:F12").Clear()
$ws2.Range("A9:F12").Clear()
